$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: shorten "Exportar los gráficos a imagen: Concentración vs Tiempo
# de creación" down to "Exportar los gráficos a imagen" (deleting the
# trailing ": Concentración vs Tiempo de creación") and leave the cursor
# right there, which is where Word drops its "_GoBack" last-edit bookmark.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(
    "Exportar los gráficos a imagen: Concentración vs Tiempo de creación",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Exportar los gráficos a imagen", 2) | Out-Null

$cursor1 = $d.Content
$cursor1.Find.Execute("Exportar los gráficos a imagen", $false, $false,
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cursor1.Collapse(0)

# Only one _GoBack bookmark can exist in the document; (re)adding it here
# moves it from its previous location (handled implicitly below because we
# rebuild that other paragraph's OOXML without the bookmark).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $cursor1) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: merge the "Poder configurar..." / "establecer un color..." runs
# (dropping the old _GoBack bookmark that used to sit between them) into one
# run, append " *Alert Values*.", and add a whole new bullet paragraph about
# implementing undo / ctrl-z.
# ---------------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("establecer un color para el gráfico.", $false, $false,
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)

$flatOpc = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00602404" w:rsidRDefault="00602404" w:rsidP="00E52059">
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr>
<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="312" w:lineRule="auto"/>
<w:jc w:val="both"/>
<w:textAlignment w:val="baseline"/>
<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t>Poder configurar los distintos gráficos de acuerdo a un valor ingresado, establecer un color para el gráfico.</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t>Alert</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t>Values</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t>.</w:t></w:r>
</w:p><w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr>
<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="312" w:lineRule="auto"/>
<w:jc w:val="both"/>
<w:textAlignment w:val="baseline"/>
<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve">Implementar la funcionalidad de </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t>undo</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t xml:space="preserve">. Agregarlo al comando rápido de </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t>ctrl</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="es-CR"/></w:rPr><w:t>-z.</w:t></w:r>
</w:p></w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($flatOpc) | Out-Null

Write-Output "edit applied"
